$d = $word.ActiveDocument
$d.Content.Find.Execute("Sample Source Code", $true, $false, $false, $false, $false, $true, 1, $false, "Примерен код", 2)
